# Apply cryptos list price/volume update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'27.935.11"
$ws.Range("E2").Value = "  +0.85%  "
# Row 3
$ws.Range("D3").Value = "'1.763.73"
$ws.Range("E3").Value = "  -0.70%  "
# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.03%  "
# Row 5
$ws.Range("D5").Value = "'328.38"
$ws.Range("E5").Value = "  +0.57%  "
# Row 6
$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  -0.07%  "
# Row 7
$ws.Range("D7").Value = "'0.4643"
$ws.Range("E7").Value = "  +0.40%  "
# Row 8
$ws.Range("D8").Value = "'0.3511"
$ws.Range("E8").Value = "  -2.13%  "
# Row 9
$ws.Range("D9").Value = "'43.87"
$ws.Range("E9").Value = "  +4.52%  "
# Row 10
$ws.Range("D10").Value = "'0.07349"
$ws.Range("E10").Value = "  -1.66%  "
# Row 11
$ws.Range("D11").Value = "'1.079"
$ws.Range("E11").Value = "  -2.04%  "
# Row 12
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.00%  "
# Row 13
$ws.Range("D13").Value = "'20.60"
$ws.Range("E13").Value = "  -0.99%  "
# Row 14
$ws.Range("D14").Value = "'5.990"
$ws.Range("E14").Value = "  -0.69%  "
# Row 15
$ws.Range("D15").Value = "'7.149"
$ws.Range("E15").Value = "  -1.28%  "
# Row 16
$ws.Range("D16").Value = "'1.762.64"
$ws.Range("E16").Value = "  -0.67%  "
# Row 17
$ws.Range("D17").Value = "'92.41"
$ws.Range("E17").Value = "  -1.27%  "
# Row 18
$ws.Range("D18").Value = "'0.00001052"
$ws.Range("E18").Value = "  -0.58%  "
# Row 19
$ws.Range("D19").Value = "'0.06417"
$ws.Range("E19").Value = "  +0.08%  "
# Row 20
$ws.Range("D20").Value = "'0.9997"
$ws.Range("E20").Value = "  -0.05%  "
# Row 21
$ws.Range("D21").Value = "'16.83"
$ws.Range("E21").Value = "  -1.51%  "
# Row 22
$ws.Range("D22").Value = "'5.757"
$ws.Range("E22").Value = "  -0.45%  "
# Row 23
$ws.Range("D23").Value = "'27.967.78"
$ws.Range("E23").Value = "  +0.68%  "
# Row 24
$ws.Range("D24").Value = "'11.12"
$ws.Range("E24").Value = "  -1.32%  "
# Row 25
$ws.Range("D25").Value = "'2.154"
$ws.Range("E25").Value = "  +3.53%  "
# Row 26
$ws.Range("D26").Value = "'162.57"
$ws.Range("E26").Value = "  -1.11%  "
# Row 27
$ws.Range("D27").Value = "'20.01"
$ws.Range("E27").Value = "  -1.65%  "
# Row 28
$ws.Range("D28").Value = "'1.966.55"
$ws.Range("E28").Value = "  -0.65%  "
# Row 29
$ws.Range("D29").Value = "'2.162"
$ws.Range("E29").Value = "  -0.21%  "
# Row 30
$ws.Range("D30").Value = "'122.84"
$ws.Range("E30").Value = "  -2.62%  "
# Row 31
$ws.Range("E31").Value = "  -2.60%  "
# Row 32
$ws.Range("D32").Value = "'0.09268"
$ws.Range("E32").Value = "  +0.46%  "
# Row 33
$ws.Range("D33").Value = "'3.647"
$ws.Range("E33").Value = "  -0.77%  "
# Row 34
$ws.Range("D34").Value = "'5.546"
$ws.Range("E34").Value = "  +0.38%  "
# Row 35
$ws.Range("D35").Value = "'11.66"
$ws.Range("E35").Value = "  -1.11%  "
# Row 36
$ws.Range("D36").Value = "'0.02267"
$ws.Range("E36").Value = "  -1.10%  "
# Row 37
$ws.Range("D37").Value = "'0.06059"
$ws.Range("E37").Value = "  -1.15%  "
# Row 38
$ws.Range("D38").Value = "'0.2059"
# Row 39
$ws.Range("D39").Value = "'4.900"
$ws.Range("E39").Value = "  -1.22%  "
# Row 40
$ws.Range("D40").Value = "'0.6121"
$ws.Range("E40").Value = "  -2.93%  "
# Row 41
$ws.Range("E41").Value = "  -0.08%  "
# Row 42
$ws.Range("B42").Value = "WEMIXTOKEN"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "'1.365"
$ws.Range("E42").Value = "  -2.02%  "
# Row 43
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'7.779"
$ws.Range("E43").Value = "  +0.09%  "
# Row 44
$ws.Range("D44").Value = "'13.14"
$ws.Range("E44").Value = "  -0.31%  "
# Row 45
$ws.Range("D45").Value = "'3.735"
$ws.Range("E45").Value = "  +0.13%  "
# Row 46
$ws.Range("D46").Value = "'0.5780"
$ws.Range("E46").Value = "  -1.79%  "
# Row 47
$ws.Range("D47").Value = "'122.75"
$ws.Range("E47").Value = "  +0.39%  "
# Row 48
$ws.Range("D48").Value = "'1.922"
$ws.Range("E48").Value = "  -1.37%  "
# Row 49
$ws.Range("D49").Value = "'0.06810"
$ws.Range("E49").Value = "  -1.73%  "
# Row 50
$ws.Range("D50").Value = "'1.120"
$ws.Range("E50").Value = "  -1.43%  "
# Row 51
$ws.Range("D51").Value = "'71.98"
$ws.Range("E51").Value = "  -0.41%  "
